# Add "90I" / "Runaway" back to the UCR offense code table.
#
# The NIBRS offense-code lookup table (sheet "UCROffenseCodeType") is
# missing the 90I "Runaway" entry. Re-insert it immediately above the
# existing 90J "Trespass of Real Property" row (i.e. as the new row 58),
# which pushes every row below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UCROffenseCodeType")

# Insert a fresh row at position 58; Excel shifts rows 58-66 down to 59-67
# and the new row inherits formatting from the row above it.
$ws.Rows.Item(58).Insert() | Out-Null

# Populate the new row with the restored 90I / Runaway entry.
$ws.Range("A58").Value = 909
$ws.Range("B58").Value = "90I"
$ws.Range("C58").Value = "Runaway"
$ws.Range("D58").Value = "Group B"
$ws.Range("E58").Value = "Other"
$ws.Range("F58").Value = "Group B Offenses (Society)"
$ws.Range("G58").Value = "Society"

# Reflect this sheet/cell as the one the author was last working in.
$ws.Activate() | Out-Null
$ws.Range("C58").Select() | Out-Null
